$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column J (2020) mirroring the style of column I for rows 3-12.
# Copy formats from I3:I12 into J3:J12 first, then set the values.
$ws.Range("I3:I12").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(4, 10).Value = 2020
$ws.Cells.Item(5, 10).Value = 253.27664777870578
$ws.Cells.Item(7, 10).Value = 93.236077839070575
$ws.Cells.Item(8, 10).Value = 160
$ws.Cells.Item(10, 10).Value = 69
$ws.Cells.Item(11, 10).Value = 48.5
$ws.Cells.Item(12, 10).Value = 22.8

# Update the sheet selection to J3 (new active cell after the edit).
$ws.Range("J3").Select()
